# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Rebuilds the worker-arrears detail table (rows 16-43) on "Hoja1" so that the
# new worker (ANA GABRIEL GAVALO EMITOLA, CC 1047388104) is interleaved with
# the existing worker (KELLY NUÑEZ AVILA, CC 45523211) period by period, and
# refreshes the "Valor Mora" / "Salario Basico" figures (columns F and G) to
# match the latest export from the EC source system.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$kelly = "45523211"
$kellyNombre = "KELLY NUÑEZ AVILA"
$ana = "1047388104"
$anaNombre = "ANA GABRIEL GAVALO EMITOLA"

# r, DocNo, Nombre, Periodo, ValorMora, SalarioBasico
$rows = @(
    @(16, $ana,   $anaNombre,   "1912", 33125, 828116),
    @(17, $ana,   $anaNombre,   "2001", 33125, 828116),
    @(18, $kelly, $kellyNombre, "2002", 35112, 877803),
    @(19, $ana,   $anaNombre,   "2002", 33125, 828116),
    @(20, $kelly, $kellyNombre, "2003", 35112, 877803),
    @(21, $ana,   $anaNombre,   "2003", 33125, 828116),
    @(22, $kelly, $kellyNombre, "2004", 35112, 877803),
    @(23, $ana,   $anaNombre,   "2004", 33125, 828116),
    @(24, $kelly, $kellyNombre, "2005", 35112, 877803),
    @(25, $ana,   $anaNombre,   "2005", 33125, 828116),
    @(26, $kelly, $kellyNombre, "2006", 35112, 877803),
    @(27, $ana,   $anaNombre,   "2006", 33125, 828116),
    @(28, $kelly, $kellyNombre, "2007", 35112, 877803),
    @(29, $ana,   $anaNombre,   "2007", 33125, 828116),
    @(30, $kelly, $kellyNombre, "2008", 35112, 877803),
    @(31, $ana,   $anaNombre,   "2008", 33125, 828116),
    @(32, $kelly, $kellyNombre, "2009", 35112, 877803),
    @(33, $ana,   $anaNombre,   "2009", 33125, 828116),
    @(34, $kelly, $kellyNombre, "2010", 35112, 877803),
    @(35, $ana,   $anaNombre,   "2010", 33125, 828116),
    @(36, $kelly, $kellyNombre, "2011", 35112, 877803),
    @(37, $ana,   $anaNombre,   "2011", 33125, 828116),
    @(38, $kelly, $kellyNombre, "2012", 35112, 877803),
    @(39, $ana,   $anaNombre,   "2012", 33125, 828116),
    @(40, $kelly, $kellyNombre, "2101", 35112, 877803),
    @(41, $ana,   $anaNombre,   "2101", 33125, 828116),
    @(42, $kelly, $kellyNombre, "2102", 25749, 877803),
    @(43, $ana,   $anaNombre,   "2102", 24292, 828116)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
}
